# Updates the "USA USL Championship" odds sheet:
#  - Several existing rows (297/298, 302/303, 309/310, 323/325) had their
#    match data (columns B..AC) swapped between the row pairs.
#  - A new match (row 327) was appended at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param(
        [int]$Row1,
        [int]$Row2
    )
    $range1 = $ws.Range("B$Row1" + ":AC$Row1")
    $range2 = $ws.Range("B$Row2" + ":AC$Row2")
    $val1 = $range1.Value()
    $val2 = $range2.Value()
    $range1.Value = $val2
    $range2.Value = $val1
}

# Swap the data between the following row pairs (row/id labels in column A
# stay put, only the match data that follows moves).
Swap-RowData 297 298
Swap-RowData 302 303
Swap-RowData 309 310
Swap-RowData 323 325

# Append the new match as row 327, copying the cell formatting (bold/border
# style for the id column, date format for the date column) from row 326.
$ws.Range("A326").Copy()
$ws.Range("A327").PasteSpecial(-4122)
$ws.Range("E326").Copy()
$ws.Range("E327").PasteSpecial(-4122)

$ws.Range("A327").Value = 325
$ws.Range("B327").Value = 7689575
$ws.Range("C327").Value = "USA USL Championship"
$ws.Range("D327").Value = "USA USL Championship"
$ws.Range("E327").Value = 45387.89583333334
$ws.Range("F327").Value = "Tulsa"
$ws.Range("G327").Value = "Phoenix Rising FC"
$ws.Range("K327").Value = 2.6
$ws.Range("L327").Value = 3.5
$ws.Range("M327").Value = 2.25
$ws.Range("N327").Value = 2.5
$ws.Range("O327").Value = 3.5
$ws.Range("P327").Value = 2.3
$ws.Range("Q327").Value = 0
$ws.Range("R327").Value = 1.975
$ws.Range("S327").Value = 1.825
$ws.Range("T327").Value = 2.75
$ws.Range("U327").Value = 1.9
$ws.Range("V327").Value = 1.9
$ws.Range("W327").Value = 0
$ws.Range("X327").Value = 0
$ws.Range("Y327").Value = 0
$ws.Range("Z327").Value = 0
$ws.Range("AA327").Value = 0

Write-Output "Edit complete"
